$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 81
$ws.Range("D81").Value = 0.7136611031493167

# Row 82
$ws.Range("D82").Value = 0.7250871391493168

# Row 83
$ws.Range("D83").Value = 0.6447426901493167

# Row 84
$ws.Range("D84").Value = 0.7064651831493167

# Row 85
$ws.Range("C85").Value = 0.5279710641493167

# Row 86
$ws.Range("C86").Value = -0.2139369238506833

# Row 87
$ws.Range("C87").Value = 0.1647551661493167

# Row 88
$ws.Range("C88").Value = 0.4518305101493167
$ws.Range("D88").Value = 0.791995474

# Row 89
$ws.Range("C89").Value = 0.4794748351493167
$ws.Range("D89").Value = 0.788120887

# Row 90
$ws.Range("C90").Value = 0.2386249091493167
$ws.Range("D90").Value = 0.597740902

# Row 91
$ws.Range("C91").Value = 0.3244906151493167
$ws.Range("D91").Value = 0.620527487

# Row 92
$ws.Range("B92").Value = 0.03316543414931669
$ws.Range("C92").Value = 0.241887844

# Row 93
$ws.Range("B93").Value = -0.0107480648506833
$ws.Range("C93").Value = 0.331651578

# Row 94
$ws.Range("C94").Value = 0.154182215

# Row 95
$ws.Range("C95").Value = 0.166899468

# Row 96
$ws.Range("C96").Value = 0.042359665

# Row 97
$ws.Range("C97").Value = 0.266698307
